$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting for the new row (row 4) from existing rows so the style
# indices match (date style for A4, plain "no-wrap" text style for B4:F4).
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A4").PasteSpecial(-4122) | Out-Null

$ws.Range("B3:C3").Copy() | Out-Null
$ws.Range("B4:C4").PasteSpecial(-4122) | Out-Null

$ws.Range("B3:C3").Copy() | Out-Null
$ws.Range("D4:F4").PasteSpecial(-4122) | Out-Null

# Fill in the new row's values.
$ws.Range("A4").Value = 43794
$ws.Range("B4").Value = "羅致遠"
$ws.Range("C4").Value = "分別控制股價走勢與買進策略等兩變因分析模擬結果"
$ws.Range("D4").Value = "完成輸出分析結果的模型"
$ws.Range("E4").Value = "協助價量模擬"
$ws.Range("F4").Value = "建立基本分析架構"

# Set the explicit height for the new row.
$ws.Rows.Item(4).RowHeight = 87

# Update selection to match the recorded cursor position after the edit.
$ws.Range("I4").Select() | Out-Null
